# Applies stock-count / value corrections produced by a later stock-take pass.
# For most items the quantity (F) was revised down slightly and the value (G)
# recomputed as Rate(D) * Qty(F); a handful of same-product rows (different item
# codes) had their batch figures (B/E/F/G) swapped between the two rows; and the
# per-company "Sub Total:" / overall "Sub Total:" / "Grand Total:" rows (column B)
# were updated to match the new sums.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 84
$ws.Range("F84").Value = 3
$ws.Range("G84").Value = 285.63

# Row 91
$ws.Range("F91").Value = 416
$ws.Range("G91").Value = 26499.2

# Row 100
$ws.Range("F100").Value = 24
$ws.Range("G100").Value = 5905.68

# Row 114
$ws.Range("B114").Value = 262982.18

# Row 136
$ws.Range("B136").Value = 63902
$ws.Range("E136").Value = 34.04
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0

# Row 137
$ws.Range("B137").Value = 48654
$ws.Range("E137").Value = 38.26
$ws.Range("F137").Value = -1
$ws.Range("G137").Value = -32.02

# Row 212
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0

# Row 222
$ws.Range("B222").Value = 50759.7

# Row 246
$ws.Range("B246").Value = 64973
$ws.Range("E246").Value = 35.4
$ws.Range("F246").Value = 64
$ws.Range("G246").Value = 2131.2

# Row 247
$ws.Range("B247").Value = 48706
$ws.Range("E247").Value = 39.8
$ws.Range("F247").Value = -144
$ws.Range("G247").Value = -4795.2

# Row 284
$ws.Range("F284").Value = 28
$ws.Range("G284").Value = 4058.04

# Row 292
$ws.Range("B292").Value = 55373
$ws.Range("E292").Value = 163.62
$ws.Range("F292").Value = -94
$ws.Range("G292").Value = -13562.32

# Row 293
$ws.Range("B293").Value = 63520
$ws.Range("E293").Value = 153.4
$ws.Range("F293").Value = 73
$ws.Range("G293").Value = 10532.44

# Row 299
$ws.Range("B299").Value = 55356
$ws.Range("E299").Value = 54.04
$ws.Range("F299").Value = -158
$ws.Range("G299").Value = -7527.12

# Row 300
$ws.Range("B300").Value = 63510
$ws.Range("E300").Value = 50.66
$ws.Range("F300").Value = 145
$ws.Range("G300").Value = 6907.8

# Row 304
$ws.Range("F304").Value = 20
$ws.Range("G304").Value = 6055.8

# Row 328
$ws.Range("F328").Value = 775
$ws.Range("G328").Value = 16298.25

# Row 330
$ws.Range("F330").Value = 5
$ws.Range("G330").Value = 2629.75

# Row 333
$ws.Range("F333").Value = 712
$ws.Range("G333").Value = 121986.96

# Row 339
$ws.Range("B339").Value = 309602.6

# Row 356
$ws.Range("B356").Value = 31930
$ws.Range("E356").Value = 26.8
$ws.Range("F356").Value = -62
$ws.Range("G356").Value = -1390.04

# Row 357
$ws.Range("B357").Value = 63681
$ws.Range("E357").Value = 23.84
$ws.Range("F357").Value = 0
$ws.Range("G357").Value = 0

# Row 382
$ws.Range("F382").Value = 156
$ws.Range("G382").Value = 6709.56

# Row 395
$ws.Range("B395").Value = 248258.1

# Row 467
$ws.Range("B467").Value = 53602
$ws.Range("E467").Value = 15.69
$ws.Range("F467").Value = -231
$ws.Range("G467").Value = -3037.65

# Row 468
$ws.Range("B468").Value = 65068
$ws.Range("E468").Value = 13.97
$ws.Range("F468").Value = 110
$ws.Range("G468").Value = 1446.5

# Row 472
$ws.Range("B472").Value = 45695
$ws.Range("E472").Value = 23.58
$ws.Range("F472").Value = -36
$ws.Range("G472").Value = -710.28

# Row 473
$ws.Range("B473").Value = 64915
$ws.Range("E473").Value = 20.98
$ws.Range("F473").Value = 0
$ws.Range("G473").Value = 0

# Row 479
$ws.Range("B479").Value = 64927
$ws.Range("E479").Value = 17.26
$ws.Range("F479").Value = 183
$ws.Range("G479").Value = 2968.26

# Row 480
$ws.Range("B480").Value = 45718
$ws.Range("E480").Value = 19.38
$ws.Range("F480").Value = -294
$ws.Range("G480").Value = -4768.68

# Row 485
$ws.Range("B485").Value = 64925
$ws.Range("E485").Value = 13.97
$ws.Range("F485").Value = 197
$ws.Range("G485").Value = 2590.55

# Row 486
$ws.Range("B486").Value = 45709
$ws.Range("E486").Value = 15.69
$ws.Range("F486").Value = -300
$ws.Range("G486").Value = -3945

# Row 490
$ws.Range("F490").Value = 245
$ws.Range("G490").Value = 3608.85

# Row 492
$ws.Range("B492").Value = -3576.6

# Row 495
$ws.Range("F495").Value = 15
$ws.Range("G495").Value = 469.5

# Row 497
$ws.Range("F497").Value = 42
$ws.Range("G497").Value = 2103.36

# Row 508
$ws.Range("B508").Value = 12810.77

# Row 559
$ws.Range("F559").Value = 102
$ws.Range("G559").Value = 2024.7

# Row 564
$ws.Range("B564").Value = 8081.49

# Row 608
$ws.Range("B608").Value = 60022
$ws.Range("E608").Value = 37.22
$ws.Range("F608").Value = -113
$ws.Range("G608").Value = -3709.79

# Row 609
$ws.Range("B609").Value = 64830
$ws.Range("E609").Value = 34.9
$ws.Range("F609").Value = 112
$ws.Range("G609").Value = 3676.96

# Row 662
$ws.Range("F662").Value = 345
$ws.Range("G662").Value = 27731.1

# Row 663
$ws.Range("B663").Value = 36443.68

# Row 710
$ws.Range("F710").Value = 41
$ws.Range("G710").Value = 3343.96

# Row 713
$ws.Range("F713").Value = 175
$ws.Range("G713").Value = 25047.75

# Row 714
$ws.Range("F714").Value = 36
$ws.Range("G714").Value = 2936.16

# Row 717
$ws.Range("B717").Value = 61428
$ws.Range("D717").Value = 69.16
$ws.Range("E717").Value = 73.52
$ws.Range("F717").Value = 1
$ws.Range("G717").Value = 69.16

# Row 718
$ws.Range("B718").Value = 63150
$ws.Range("D718").Value = 75.68000000000001
$ws.Range("E718").Value = 80.45
$ws.Range("F718").Value = 64
$ws.Range("G718").Value = 4843.52

# Row 727
$ws.Range("F727").Value = 264
$ws.Range("G727").Value = 31867.44

# Row 728
$ws.Range("B728").Value = 140152.56

# Row 732
$ws.Range("F732").Value = 45
$ws.Range("G732").Value = 7375.05

# Row 754
$ws.Range("F754").Value = 112
$ws.Range("G754").Value = 6382.88

# Row 755
$ws.Range("B755").Value = 80932.06

# Row 780
$ws.Range("F780").Value = 3169
$ws.Range("G780").Value = 516895.59

# Row 782
$ws.Range("F782").Value = 593
$ws.Range("G782").Value = 167741.91

# Row 787
$ws.Range("B787").Value = 787599.7

# Row 790
$ws.Range("F790").Value = 102
$ws.Range("G790").Value = 14893.02

# Row 804
$ws.Range("B804").Value = 78537.10000000001

# Row 805
$ws.Range("B805").Value = 3057548.72

# Row 806
$ws.Range("B806").Value = 3057548.72
